$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Insert a new "Meta description" paragraph right after the title
#    (first paragraph, Heading1 style).
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
[void]$titlePara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs(2)

$metaXml = '<?xml version="1.0" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p>' +
    '<w:r/>' +
    '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' +
    '<w:r><w:t>: Read our review of the online slot game Beat the Beast: Cerberus'' Inferno, and play for free with exciting features, great payouts, and a Greek mythology theme.</w:t></w:r>' +
    '</w:p>' +
    '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

[void]$newPara.Range.InsertXML($metaXml)

# ------------------------------------------------------------------
# 2. Remove the duplicated bold title paragraph near the end of the
#    document, and replace the text of the following (italic)
#    paragraph with the new image-prompt text.
# ------------------------------------------------------------------
$cr = [char]13
$titleText = "Play Beat the Beast: Cerberus' Inferno for Free - Review"
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text.TrimEnd($cr)
    if ($t -eq $titleText -and $i -gt 1) {
        [void]$p.Range.Delete()
        break
    }
}

$oldDescription = "Read our review of the online slot game Beat the Beast: Cerberus' Inferno, and play for free with exciting features, great payouts, and a Greek mythology theme."
$newDescription = "Create an eye-catching feature image in cartoon style that features a happy Maya warrior with glasses, armed with a sword and shield, standing triumphantly in front of a fiery inferno backdrop. The warrior should be wearing a headband with the game title ""Beat the Beast: Cerberus' Inferno"" written on it. The overall look should be lively and vibrant, with bright colors and bold lines to capture the adventurous feel of the game."

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text.TrimEnd($cr)
    if ($t -eq $oldDescription) {
        $rng = $p.Range
        $target = $d.Range($rng.Start, $rng.End - 1)
        $target.Text = $newDescription
        break
    }
}

Write-Output "done"
